$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.923.09'
$ws.Range('E2').Value = '  -4.00%  '
$ws.Range('D3').Value = '3.514.63'
$ws.Range('E3').Value = '  -4.90%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '579.82'
$ws.Range('E5').Value = '  -1.61%  '
$ws.Range('D6').Value = '174.09'
$ws.Range('E6').Value = '  -2.84%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '3.506.73'
$ws.Range('E8').Value = '  -4.93%  '
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('E10').Value = '  -6.54%  '
$ws.Range('E11').Value = '  +7.81%  '
$ws.Range('E12').Value = '  -2.54%  '
$ws.Range('D13').Value = '47.25'
$ws.Range('E13').Value = '  -5.51%  '
$ws.Range('E14').Value = '  -3.54%  '
$ws.Range('D15').Value = '671.15'
$ws.Range('E15').Value = '  -1.54%  '
$ws.Range('D16').Value = '4.080.20'
$ws.Range('E16').Value = '  -4.78%  '
$ws.Range('D17').Value = '8.80'
$ws.Range('E17').Value = '  -2.13%  '
$ws.Range('D18').Value = '3.513.35'
$ws.Range('E18').Value = '  -4.14%  '
$ws.Range('D19').Value = '68.952.39'
$ws.Range('E19').Value = '  -4.12%  '
$ws.Range('E20').Value = '  -1.53%  '
$ws.Range('D21').Value = '17.54'
$ws.Range('E21').Value = '  -2.93%  '
$ws.Range('E22').Value = '  -3.67%  '
$ws.Range('D23').Value = '0.906'
$ws.Range('E23').Value = '  -3.56%  '
$ws.Range('D24').Value = '16.28'
$ws.Range('E24').Value = '  -8.66%  '
$ws.Range('D25').Value = '98.39'
$ws.Range('E25').Value = '  -4.93%  '
$ws.Range('E26').Value = '  -4.37%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('E28').Value = '  -6.59%  '
$ws.Range('D29').Value = '9.44'
$ws.Range('E29').Value = '  -7.27%  '
$ws.Range('D30').Value = '32.96'
$ws.Range('E30').Value = '  -7.23%  '
$ws.Range('E31').Value = '  -4.81%  '
$ws.Range('E32').Value = '  -7.72%  '
$ws.Range('E34').Value = '  -4.85%  '
$ws.Range('D35').Value = '577.98'
$ws.Range('E35').Value = '  +0.41%  '
$ws.Range('E36').Value = '  -3.31%  '
$ws.Range('E37').Value = '  -14.52%  '
$ws.Range('E38').Value = '  -4.17%  '
$ws.Range('D39').Value = '57.15'
$ws.Range('E39').Value = '  -4.09%  '
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.17%  '
$ws.Range('E41').Value = '  -3.55%  '
$ws.Range('E42').Value = '  -5.42%  '
$ws.Range('E43').Value = '  -6.21%  '
$ws.Range('D44').Value = '3.418.07'
$ws.Range('E44').Value = '  -9.23%  '
$ws.Range('D45').Value = '33.47'
$ws.Range('E45').Value = '  -5.62%  '
$ws.Range('D46').Value = '0.0₃0704'
$ws.Range('E46').Value = '  -9.38%  '
$ws.Range('D47').Value = '2.90'
$ws.Range('E47').Value = '  +0.29%  '
$ws.Range('E48').Value = '  -7.09%  '
$ws.Range('E49').Value = '  -0.43%  '
$ws.Range('D50').Value = '131.53'
$ws.Range('E50').Value = '  -2.08%  '
$ws.Range('E51').Value = '  -0.02%  '
